$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the slightly adjusted timestamp on row 5 (floating point precision refresh)
$ws.Range("A5").Value = 45864.37523967592

# Append new row 6 with the latest sensor reading
$ws.Range("A6").Value = 45864.41694890593
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 13.89
$ws.Range("E6").Value = 89.01000000000001
$ws.Range("F6").Value = 520.3099999999999
$ws.Range("G6").Value = 3.33
$ws.Range("H6").Value = "ESE"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "10:00:24"

# Match the number format used by the other timestamp cells in column A
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
